$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1969111969111969
$ws.Range("C2").Value = 0.5366795366795367
$ws.Range("J2").Value = 0.003861003861003861
$ws.Range("P2").Value = 0.1814671814671815
$ws.Range("S2").Value = 0.08108108108108109
$ws.Range("B3").Value = 0.01388888888888889
$ws.Range("C3").Value = 0.02777777777777778
$ws.Range("J3").Value = 0.02083333333333333
$ws.Range("P3").Value = 0.7152777777777778
$ws.Range("S3").Value = 0.2222222222222222
$ws.Range("J4").Value = 0.02040816326530612
$ws.Range("O4").Value = 0.02040816326530612
$ws.Range("P4").Value = 0.6530612244897959
$ws.Range("S4").Value = 0.3061224489795918
$ws.Range("B6").Value = 0.064
$ws.Range("D6").Value = 0.012
$ws.Range("F6").Value = 0.104
$ws.Range("J6").Value = 0.144
$ws.Range("O6").Value = 0.02
$ws.Range("Q6").Value = 0.156
$ws.Range("R6").Value = 0.116
$ws.Range("S6").Value = 0.384
$ws.Range("B7").Value = 0.125
$ws.Range("F7").Value = 0.08125
$ws.Range("J7").Value = 0.13125
$ws.Range("O7").Value = 0.01875
$ws.Range("R7").Value = 0.08749999999999999
$ws.Range("S7").Value = 0.43125
$ws.Range("B8").Value = 0.07142857142857142
$ws.Range("D8").Value = 0.01724137931034483
$ws.Range("F8").Value = 0.06896551724137931
$ws.Range("J8").Value = 0.145320197044335
$ws.Range("O8").Value = 0.02463054187192118
$ws.Range("Q8").Value = 0.1650246305418719
$ws.Range("R8").Value = 0.1133004926108374
$ws.Range("S8").Value = 0.3940886699507389
$ws.Range("B9").Value = 0.07630522088353414
$ws.Range("D9").Value = 0.02409638554216868
$ws.Range("F9").Value = 0.1124497991967871
$ws.Range("J9").Value = 0.1044176706827309
$ws.Range("O9").Value = 0.0321285140562249
$ws.Range("Q9").Value = 0.2008032128514056
$ws.Range("R9").Value = 0.09236947791164658
$ws.Range("S9").Value = 0.357429718875502
$ws.Range("B10").Value = 0.09304207119741101
$ws.Range("D10").Value = 0.02427184466019417
$ws.Range("E10").Value = 0.0008090614886731392
$ws.Range("F10").Value = 0.09385113268608414
$ws.Range("J10").Value = 0.1059870550161812
$ws.Range("O10").Value = 0.01537216828478964
$ws.Range("Q10").Value = 0.1868932038834951
$ws.Range("R10").Value = 0.1108414239482201
$ws.Range("S10").Value = 0.3689320388349515
$ws.Range("G11").Value = 0.1363636363636364
$ws.Range("J11").Value = 0.1136363636363636
$ws.Range("K11").Value = 0.1931818181818182
$ws.Range("L11").Value = 0.5378787878787878
$ws.Range("S11").Value = 0.01893939393939394
$ws.Range("G12").Value = 0.6754966887417219
$ws.Range("J12").Value = 0.2450331125827815
$ws.Range("K12").Value = 0.006622516556291391
$ws.Range("L12").Value = 0.05960264900662252
$ws.Range("S12").Value = 0.01324503311258278
$ws.Range("F13").Value = 0.02325581395348837
$ws.Range("G13").Value = 0.6744186046511628
$ws.Range("J13").Value = 0.2558139534883721
$ws.Range("S13").Value = 0.04651162790697674
$ws.Range("F15").Value = 0.01476014760147601
$ws.Range("H15").Value = 0.1143911439114391
$ws.Range("I15").Value = 0.1033210332103321
$ws.Range("J15").Value = 0.3874538745387454
$ws.Range("K15").Value = 0.04797047970479705
$ws.Range("M15").Value = 0.007380073800738007
$ws.Range("O15").Value = 0.1070110701107011
$ws.Range("S15").Value = 0.2177121771217712
$ws.Range("F16").Value = 0.01129943502824859
$ws.Range("H16").Value = 0.1581920903954802
$ws.Range("I16").Value = 0.1412429378531073
$ws.Range("J16").Value = 0.3615819209039548
$ws.Range("K16").Value = 0.1073446327683616
$ws.Range("M16").Value = 0.02259887005649718
$ws.Range("N16").Value = 0.005649717514124294
$ws.Range("O16").Value = 0.06214689265536723
$ws.Range("S16").Value = 0.1299435028248588
$ws.Range("F17").Value = 0.0125
$ws.Range("H17").Value = 0.1525
$ws.Range("I17").Value = 0.1125
$ws.Range("J17").Value = 0.44
$ws.Range("K17").Value = 0.1075
$ws.Range("M17").Value = 0.02
$ws.Range("N17").Value = 0.0025
$ws.Range("O17").Value = 0.05
$ws.Range("S17").Value = 0.1025
$ws.Range("F18").Value = 0.02834008097165992
$ws.Range("H18").Value = 0.1700404858299595
$ws.Range("I18").Value = 0.09716599190283401
$ws.Range("J18").Value = 0.4696356275303644
$ws.Range("K18").Value = 0.06477732793522267
$ws.Range("M18").Value = 0.01619433198380567
$ws.Range("N18").Value = 0.004048582995951417
$ws.Range("O18").Value = 0.0728744939271255
$ws.Range("S18").Value = 0.07692307692307693
$ws.Range("F19").Value = 0.006739679865206402
$ws.Range("H19").Value = 0.2055602358887953
$ws.Range("I19").Value = 0.1069924178601516
$ws.Range("J19").Value = 0.3681550126368998
$ws.Range("K19").Value = 0.101095197978096
$ws.Range("M19").Value = 0.02190395956192081
$ws.Range("N19").Value = 0.002527379949452401
$ws.Range("O19").Value = 0.09267059814658804
$ws.Range("S19").Value = 0.09435551811288964
